# Add a new "Intermediate SQL Queries" course column (G) to the datacamp
# course-ratings sheet, mirroring the existing layout:
#   - new header in G1 (adds a new shared string)
#   - new rating value in G2
#   - column G sized/auto-fit like the other bestFit columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Intermediate SQL Queries"
$ws.Range("G2").Value = 4

# Size column G the way Excel would after a best-fit autofit on the new
# header text (matches the style of columns A:F, which are all bestFit).
$ws.Columns.Item(7).ColumnWidth = 23.5
